{"js": "// Corrections to both files:\n//  1. Remove the stray \"_GoBack\" bookmark that sat after \"WORK EXPERIENCE\".\n//  2. Add a comma right after \"IAVVC\" (so the citation reads \"IAVVC, 2023\")\n//     and re-insert the \"_GoBack\" bookmark at that new location (directly\n//     after the inserted comma, before the following space/\"2023\").\n\nconst doc = context.document;\n\n// 1) Drop the old bookmark wherever it currently lives.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the \"IAVVC\" run so we can insert the missing comma right after it.\nconst searchResults = doc.body.search(\"IAVVC\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const iavvcRange = searchResults.items[0];\n\n  // Insert \",\" immediately after \"IAVVC\"; insertText returns the range of\n  // the newly inserted text.\n  const commaRange = iavvcRange.insertText(\",\", Word.InsertLocation.after);\n  await context.sync();\n\n  // Re-create the \"_GoBack\" bookmark collapsed right after the new comma.\n  const afterComma = commaRange.getRange(Word.RangeLocation.end);\n  afterComma.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Corrections to both files:\n#  1. Remove the stray \"_GoBack\" bookmark that sat after \"WORK EXPERIENCE\".\n#  2. Add a comma right after \"IAVVC\" (so the citation reads \"IAVVC, 2023\")\n#     and re-insert the \"_GoBack\" bookmark at that new location (directly\n#     after the inserted comma, before the following space/\"2023\").\n\n$d = $word.ActiveDocument\n\n# 1) Drop the old bookmark wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Find \"IAVVC\" and insert the missing comma right after it.\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"IAVVC\"\n$find.Execute() | Out-Null\n\n$insertRange = $d.Range($findRange.End, $findRange.End)\n$insertRange.InsertAfter(\",\")\n\n# 3) Re-locate \"IAVVC,\" (now that the comma is in the document) so we get a\n#    fresh, reliable range to anchor the bookmark on.\n$bmFindRange = $d.Content\n$bmFind = $bmFindRange.Find\n$bmFind.Text = \"IAVVC,\"\n$bmFind.Execute() | Out-Null\n\n$bookmarkRange = $d.Range($bmFindRange.End, $bmFindRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
